# Atualizacao rapida de agenda as  9:13:02,55
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Giovani / Residencial Santos Drumont
$ws.Range("A2").Value = "Giovani"
$ws.Range("B2").Value = "'0884"
$ws.Range("C2").Value = "Residencial Santos Drumont"
$ws.Range("D2").Value = "Sem comunicação geral, não sabemos dizer se é internet."
$ws.Range("E2").Value = ""
$ws.Range("G2").Value = "Pendente"
$ws.Range("H2").Value = "Maxvel: 19 / Forte: 16"

# Row 3 - Giovani / Igreja Imaculada
$ws.Range("A3").Value = "Giovani"
$ws.Range("B3").Value = "'0390"
$ws.Range("C3").Value = "Igreja Imaculada"
$ws.Range("D3").Value = "Sem comunicação de alarmes, tentar atualizar a central."
$ws.Range("E3").Value = ""
$ws.Range("G3").Value = "Pendente"
$ws.Range("H3").Value = ""

# Row 4 - Giovani / Valinhos Departamento Pessoal
$ws.Range("A4").Value = "Giovani"
$ws.Range("B4").Value = "'0554"
$ws.Range("C4").Value = "Valinhos Departamento Pessoal"
$ws.Range("D4").Value = "Sem comunicação de alarmes, via linha telefônica."
$ws.Range("E4").Value = ""
$ws.Range("G4").Value = "Pendente"
$ws.Range("H4").Value = ""

# Row 5 - Giovani / Viaceu Loja
$ws.Range("A5").Value = "Giovani"
$ws.Range("B5").Value = "'0288"
$ws.Range("C5").Value = "Viaceu Loja"
$ws.Range("D5").Value = "Sem comunicação de câmeras."
$ws.Range("E5").Value = ""
$ws.Range("G5").Value = "Pendente"
$ws.Range("H5").Value = ""

# Row 6 - Giovani / Brapi
$ws.Range("A6").Value = "Giovani"
$ws.Range("B6").Value = "'0217"
$ws.Range("C6").Value = "Brapi"
$ws.Range("D6").Value = "Zonas abertas, cliente pedindo reparo."
$ws.Range("E6").Value = ""
$ws.Range("G6").Value = "Pendente"
$ws.Range("H6").Value = ""

# Row 7 - Roberto / Galpão JMC
$ws.Range("A7").Value = "Roberto"
$ws.Range("B7").Value = "'0786"
$ws.Range("C7").Value = "Galpão JMC"
$ws.Range("D7").Value = "Cliente pedindo reparo, controles sem funcionar, zonas abertas, sensor quebrado, cabo arrebentado, DVR travado e pode te rmais detalhes."
$ws.Range("E7").Value = ""
$ws.Range("G7").Value = "Pendente"
$ws.Range("H7").Value = ""

# Row 8 - Marcos / Sitio Alves
$ws.Range("A8").Value = "Marcos"
$ws.Range("B8").Value = "'0422"
$ws.Range("C8").Value = "Sitio Alves"
$ws.Range("D8").Value = "Sem comunicação de câmeras, zona aberta e central no GPRS."
$ws.Range("E8").Value = "Imagens restauradas, central via internet novamente."
$ws.Range("G8").Value = "Pendente"
$ws.Range("H8").Value = ""

# Rows 9-11: clear entirely (fewer agenda entries now)
$ws.Range("A9:H11").ClearContents()

# Update sheet view scroll/selection to match the edited session
$ws.Range("A9").Select()
